$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates derived from the crypto-price refresh diff.
# Force Text number format before writing so Excel does not
# auto-convert numeric-looking strings (e.g. "26.944.59", "1.001")
# into real numbers, then restore the cell's original (Normal) style
# so no stray formatting is left behind.
$updates = @(
    @{Cell='D2'; Value='26.944.59'},
    @{Cell='E2'; Value='  -3.60%  '},
    @{Cell='D3'; Value='1.715.28'},
    @{Cell='E3'; Value='  -3.01%  '},
    @{Cell='D4'; Value='1.001'},
    @{Cell='E4'; Value='  -0.02%  '},
    @{Cell='D5'; Value='307.86'},
    @{Cell='E5'; Value='  -6.41%  '},
    @{Cell='E6'; Value='  +0.06%  '},
    @{Cell='D7'; Value='0.4811'},
    @{Cell='E7'; Value='  +5.26%  '},
    @{Cell='D8'; Value='0.3477'},
    @{Cell='E8'; Value='  -1.37%  '},
    @{Cell='D9'; Value='41.93'},
    @{Cell='D10'; Value='0.07241'},
    @{Cell='E10'; Value='  -2.06%  '},
    @{Cell='D11'; Value='1.041'},
    @{Cell='E11'; Value='  -4.98%  '},
    @{Cell='E12'; Value='  +0.02%  '},
    @{Cell='E13'; Value='  -4.39%  '},
    @{Cell='D14'; Value='5.833'},
    @{Cell='E14'; Value='  -2.91%  '},
    @{Cell='D15'; Value='1.714.41'},
    @{Cell='E15'; Value='  -2.82%  '},
    @{Cell='D16'; Value='6.834'},
    @{Cell='E16'; Value='  -4.99%  '},
    @{Cell='E17'; Value='  -2.36%  '},
    @{Cell='D18'; Value='86.12'},
    @{Cell='E18'; Value='  -7.09%  '},
    @{Cell='E19'; Value='  -1.18%  '},
    @{Cell='E20'; Value='  +0.03%  '},
    @{Cell='D21'; Value='16.49'},
    @{Cell='E21'; Value='  -2.51%  '},
    @{Cell='D22'; Value='5.613'},
    @{Cell='E22'; Value='  -2.79%  '},
    @{Cell='D23'; Value='27.008.81'},
    @{Cell='E23'; Value='  -3.44%  '},
    @{Cell='D24'; Value='10.74'},
    @{Cell='E24'; Value='  -4.28%  '},
    @{Cell='D25'; Value='2.087'},
    @{Cell='E25'; Value='  -2.41%  '},
    @{Cell='D26'; Value='152.28'},
    @{Cell='E26'; Value='  -5.61%  '},
    @{Cell='D27'; Value='19.94'},
    @{Cell='E27'; Value='  -1.17%  '},
    @{Cell='D28'; Value='1.910.58'},
    @{Cell='E28'; Value='  -2.97%  '},
    @{Cell='D29'; Value='2.073'},
    @{Cell='E29'; Value='  -4.21%  '},
    @{Cell='D30'; Value='120.99'},
    @{Cell='E30'; Value='  -2.54%  '},
    @{Cell='D31'; Value='1.024'},
    @{Cell='E31'; Value='  -4.85%  '},
    @{Cell='D32'; Value='0.09152'},
    @{Cell='E32'; Value='  -1.53%  '},
    @{Cell='D33'; Value='3.594'},
    @{Cell='E33'; Value='  -2.00%  '},
    @{Cell='D34'; Value='5.319'},
    @{Cell='E34'; Value='  -5.06%  '},
    @{Cell='D35'; Value='1.468'},
    @{Cell='E35'; Value='  +6.48%  '},
    @{Cell='D36'; Value='0.02174'},
    @{Cell='E36'; Value='  -4.66%  '},
    @{Cell='D37'; Value='0.05858'},
    @{Cell='E37'; Value='  -4.49%  '},
    @{Cell='D38'; Value='0.2002'},
    @{Cell='E38'; Value='  -4.02%  '},
    @{Cell='B39'; Value='Aptos'},
    @{Cell='C39'; Value='https://coinranking.com/coin/HGYj5JCv5+aptos-apt'},
    @{Cell='D39'; Value='10.93'},
    @{Cell='E39'; Value='  -7.86%  '},
    @{Cell='B40'; Value='TheSandbox'},
    @{Cell='C40'; Value='https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'},
    @{Cell='D40'; Value='0.6044'},
    @{Cell='E40'; Value='  -3.48%  '},
    @{Cell='D41'; Value='4.720'},
    @{Cell='E41'; Value='  -4.41%  '},
    @{Cell='E42'; Value='  -8.04%  '},
    @{Cell='D43'; Value='7.428'},
    @{Cell='E43'; Value='  -4.96%  '},
    @{Cell='D44'; Value='12.74'},
    @{Cell='E44'; Value='  -3.13%  '},
    @{Cell='E45'; Value='  -4.71%  '},
    @{Cell='D46'; Value='0.5630'},
    @{Cell='E46'; Value='  -3.86%  '},
    @{Cell='D47'; Value='118.75'},
    @{Cell='E47'; Value='  -2.97%  '},
    @{Cell='E48'; Value='  -5.69%  '},
    @{Cell='E49'; Value='  -2.16%  '},
    @{Cell='D50'; Value='0.06653'},
    @{Cell='E50'; Value='  -2.38%  '},
    @{Cell='E51'; Value='  +0.15%  '}
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $u.Value
    $rng.Style = "Normal"
}
